# Apply 2025-06-06 data update to violent-crime-full-year.xlsx
# Generated from the unified OOXML diff; each entry updates one cell
# in the "2025" (YTD) column (L) or a revised "2024" (K) value.

$wb = $excel.ActiveWorkbook

$changes = @(
    @{ Sheet = 'Citywide Totals'; Cell = 'L2'; Value = 2713 }
    @{ Sheet = 'Citywide Totals'; Cell = 'L3'; Value = 2750 }
    @{ Sheet = 'Citywide Totals'; Cell = 'K4'; Value = 1766 }
    @{ Sheet = 'Citywide Totals'; Cell = 'L4'; Value = 733 }
    @{ Sheet = 'Citywide Totals'; Cell = 'L6'; Value = 2462 }
    @{ Sheet = 'Citywide Totals'; Cell = 'K7'; Value = 27558 }
    @{ Sheet = 'Citywide Totals'; Cell = 'L7'; Value = 8818 }
    @{ Sheet = 'By Neighborhood'; Cell = 'L2'; Value = 71 }
    @{ Sheet = 'By Neighborhood'; Cell = 'L6'; Value = 68 }
    @{ Sheet = 'By Neighborhood'; Cell = 'L7'; Value = 291 }
    @{ Sheet = 'By Neighborhood'; Cell = 'L8'; Value = 558 }
    @{ Sheet = 'By Neighborhood'; Cell = 'L11'; Value = 154 }
    @{ Sheet = 'By Neighborhood'; Cell = 'L19'; Value = 247 }
    @{ Sheet = 'By Neighborhood'; Cell = 'L20'; Value = 223 }
    @{ Sheet = 'By Neighborhood'; Cell = 'L29'; Value = 471 }
    @{ Sheet = 'By Neighborhood'; Cell = 'L30'; Value = 42 }
    @{ Sheet = 'By Neighborhood'; Cell = 'L33'; Value = 398 }
    @{ Sheet = 'By Neighborhood'; Cell = 'L35'; Value = 13 }
    @{ Sheet = 'By Neighborhood'; Cell = 'L36'; Value = 122 }
    @{ Sheet = 'By Neighborhood'; Cell = 'L37'; Value = 326 }
    @{ Sheet = 'By Neighborhood'; Cell = 'L42'; Value = 291 }
    @{ Sheet = 'By Neighborhood'; Cell = 'L51'; Value = 106 }
    @{ Sheet = 'By Neighborhood'; Cell = 'L52'; Value = 178 }
    @{ Sheet = 'By Neighborhood'; Cell = 'L53'; Value = 105 }
    @{ Sheet = 'By Neighborhood'; Cell = 'L54'; Value = 180 }
    @{ Sheet = 'By Neighborhood'; Cell = 'L55'; Value = 85 }
    @{ Sheet = 'By Neighborhood'; Cell = 'L57'; Value = 35 }
    @{ Sheet = 'By Neighborhood'; Cell = 'L60'; Value = 53 }
    @{ Sheet = 'By Neighborhood'; Cell = 'K63'; Value = 159 }
    @{ Sheet = 'By Neighborhood'; Cell = 'L63'; Value = 28 }
    @{ Sheet = 'By Neighborhood'; Cell = 'L65'; Value = 162 }
    @{ Sheet = 'By Neighborhood'; Cell = 'L67'; Value = 326 }
    @{ Sheet = 'By Neighborhood'; Cell = 'L71'; Value = 24 }
    @{ Sheet = 'By Neighborhood'; Cell = 'L72'; Value = 38 }
    @{ Sheet = 'By Neighborhood'; Cell = 'L77'; Value = 53 }
    @{ Sheet = 'By Neighborhood'; Cell = 'L78'; Value = 115 }
    @{ Sheet = 'By Neighborhood'; Cell = 'L79'; Value = 237 }
    @{ Sheet = 'By Neighborhood'; Cell = 'L83'; Value = 209 }
    @{ Sheet = 'By Neighborhood'; Cell = 'L84'; Value = 93 }
    @{ Sheet = 'By Neighborhood'; Cell = 'L85'; Value = 458 }
    @{ Sheet = 'By Neighborhood'; Cell = 'L89'; Value = 113 }
    @{ Sheet = 'By Neighborhood'; Cell = 'L90'; Value = 86 }
    @{ Sheet = 'By Neighborhood'; Cell = 'L91'; Value = 125 }
    @{ Sheet = 'By Neighborhood'; Cell = 'L92'; Value = 25 }
    @{ Sheet = 'By Neighborhood'; Cell = 'L99'; Value = 146 }
    @{ Sheet = 'By Neighborhood'; Cell = 'K101'; Value = 27558 }
    @{ Sheet = 'By Neighborhood'; Cell = 'L101'; Value = 8818 }
    @{ Sheet = 'Auburn Gresham'; Cell = 'L2'; Value = 87 }
    @{ Sheet = 'Auburn Gresham'; Cell = 'L3'; Value = 93 }
    @{ Sheet = 'Auburn Gresham'; Cell = 'L7'; Value = 291 }
    @{ Sheet = 'Belmont Cragin'; Cell = 'L3'; Value = 49 }
    @{ Sheet = 'Belmont Cragin'; Cell = 'L7'; Value = 154 }
    @{ Sheet = 'Uptown'; Cell = 'L3'; Value = 29 }
    @{ Sheet = 'Uptown'; Cell = 'L7'; Value = 113 }
    @{ Sheet = 'South Shore'; Cell = 'L3'; Value = 186 }
    @{ Sheet = 'South Shore'; Cell = 'L7'; Value = 458 }
    @{ Sheet = 'Little Village'; Cell = 'L2'; Value = 60 }
    @{ Sheet = 'Little Village'; Cell = 'L6'; Value = 49 }
    @{ Sheet = 'Little Village'; Cell = 'L7'; Value = 178 }
    @{ Sheet = 'Logan Square'; Cell = 'L3'; Value = 26 }
    @{ Sheet = 'Logan Square'; Cell = 'L6'; Value = 35 }
    @{ Sheet = 'Logan Square'; Cell = 'L7'; Value = 105 }
    @{ Sheet = 'Austin'; Cell = 'L2'; Value = 160 }
    @{ Sheet = 'Austin'; Cell = 'L7'; Value = 558 }
    @{ Sheet = 'South Chicago'; Cell = 'L3'; Value = 88 }
    @{ Sheet = 'South Chicago'; Cell = 'L6'; Value = 45 }
    @{ Sheet = 'South Chicago'; Cell = 'L7'; Value = 209 }
    @{ Sheet = 'Garfield Park'; Cell = 'L2'; Value = 108 }
    @{ Sheet = 'Garfield Park'; Cell = 'L6'; Value = 136 }
    @{ Sheet = 'Garfield Park'; Cell = 'L7'; Value = 398 }
    @{ Sheet = 'Grand Crossing'; Cell = 'L2'; Value = 96 }
    @{ Sheet = 'Grand Crossing'; Cell = 'L6'; Value = 105 }
    @{ Sheet = 'Grand Crossing'; Cell = 'L7'; Value = 326 }
    @{ Sheet = 'New City'; Cell = 'L3'; Value = 51 }
    @{ Sheet = 'New City'; Cell = 'L4'; Value = 7 }
    @{ Sheet = 'New City'; Cell = 'L7'; Value = 162 }
    @{ Sheet = 'Woodlawn'; Cell = 'L4'; Value = 14 }
    @{ Sheet = 'Woodlawn'; Cell = 'L7'; Value = 146 }
    @{ Sheet = 'Fuller Park'; Cell = 'L2'; Value = 15 }
    @{ Sheet = 'Fuller Park'; Cell = 'L7'; Value = 42 }
    @{ Sheet = 'North Lawndale'; Cell = 'L2'; Value = 95 }
    @{ Sheet = 'North Lawndale'; Cell = 'L3'; Value = 117 }
    @{ Sheet = 'North Lawndale'; Cell = 'L7'; Value = 326 }
    @{ Sheet = 'South Deering'; Cell = 'L3'; Value = 37 }
    @{ Sheet = 'South Deering'; Cell = 'L7'; Value = 93 }
    @{ Sheet = 'Loop'; Cell = 'L3'; Value = 35 }
    @{ Sheet = 'Loop'; Cell = 'L4'; Value = 14 }
    @{ Sheet = 'Loop'; Cell = 'L6'; Value = 92 }
    @{ Sheet = 'Loop'; Cell = 'L7'; Value = 180 }
    @{ Sheet = 'Englewood'; Cell = 'L2'; Value = 150 }
    @{ Sheet = 'Englewood'; Cell = 'L3'; Value = 175 }
    @{ Sheet = 'Englewood'; Cell = 'L6'; Value = 120 }
    @{ Sheet = 'Englewood'; Cell = 'L7'; Value = 471 }
    @{ Sheet = 'Chatham'; Cell = 'L2'; Value = 83 }
    @{ Sheet = 'Chatham'; Cell = 'L3'; Value = 77 }
    @{ Sheet = 'Chatham'; Cell = 'L7'; Value = 247 }
    @{ Sheet = 'Ashburn'; Cell = 'L3'; Value = 20 }
    @{ Sheet = 'Ashburn'; Cell = 'L6'; Value = 14 }
    @{ Sheet = 'Ashburn'; Cell = 'L7'; Value = 68 }
    @{ Sheet = 'Humboldt Park'; Cell = 'L2'; Value = 83 }
    @{ Sheet = 'Humboldt Park'; Cell = 'L6'; Value = 86 }
    @{ Sheet = 'Humboldt Park'; Cell = 'L7'; Value = 291 }
    @{ Sheet = 'Rogers Park'; Cell = 'L3'; Value = 33 }
    @{ Sheet = 'Rogers Park'; Cell = 'L7'; Value = 115 }
    @{ Sheet = 'Lower West Side'; Cell = 'L2'; Value = 32 }
    @{ Sheet = 'Lower West Side'; Cell = 'L7'; Value = 85 }
    @{ Sheet = 'Washington Park'; Cell = 'L3'; Value = 47 }
    @{ Sheet = 'Washington Park'; Cell = 'L7'; Value = 125 }
    @{ Sheet = 'Roseland'; Cell = 'L2'; Value = 77 }
    @{ Sheet = 'Roseland'; Cell = 'L6'; Value = 47 }
    @{ Sheet = 'Roseland'; Cell = 'L7'; Value = 237 }
    @{ Sheet = 'Chicago Lawn'; Cell = 'L3'; Value = 69 }
    @{ Sheet = 'Chicago Lawn'; Cell = 'L7'; Value = 223 }
    @{ Sheet = 'Grand Boulevard'; Cell = 'L3'; Value = 30 }
    @{ Sheet = 'Grand Boulevard'; Cell = 'L6'; Value = 31 }
    @{ Sheet = 'Grand Boulevard'; Cell = 'L7'; Value = 122 }
    @{ Sheet = 'Gold Coast'; Cell = 'L4'; Value = 6 }
    @{ Sheet = 'Gold Coast'; Cell = 'L7'; Value = 13 }
    @{ Sheet = 'Albany Park'; Cell = 'L3'; Value = 23 }
    @{ Sheet = 'Albany Park'; Cell = 'L7'; Value = 71 }
    @{ Sheet = 'West Elsdon'; Cell = 'L2'; Value = 9 }
    @{ Sheet = 'West Elsdon'; Cell = 'L7'; Value = 25 }
    @{ Sheet = 'Washington Heights'; Cell = 'L6'; Value = 23 }
    @{ Sheet = 'Washington Heights'; Cell = 'L7'; Value = 86 }
    @{ Sheet = 'Little Italy, UIC'; Cell = 'L6'; Value = 28 }
    @{ Sheet = 'Little Italy, UIC'; Cell = 'L7'; Value = 106 }
    @{ Sheet = 'Mckinley Park'; Cell = 'L3'; Value = 9 }
    @{ Sheet = 'Mckinley Park'; Cell = 'L7'; Value = 35 }
    @{ Sheet = 'Morgan Park'; Cell = 'L3'; Value = 19 }
    @{ Sheet = 'Morgan Park'; Cell = 'L7'; Value = 53 }
    @{ Sheet = 'Oakland'; Cell = 'L2'; Value = 8 }
    @{ Sheet = 'Oakland'; Cell = 'L7'; Value = 24 }
    @{ Sheet = 'Old Town'; Cell = 'L4'; Value = 6 }
    @{ Sheet = 'Old Town'; Cell = 'L6'; Value = 10 }
    @{ Sheet = 'Old Town'; Cell = 'L7'; Value = 38 }
    @{ Sheet = 'Riverdale'; Cell = 'L3'; Value = 19 }
    @{ Sheet = 'Riverdale'; Cell = 'L7'; Value = 53 }
)

foreach ($chg in $changes) {
    $ws = $wb.Worksheets($chg.Sheet)
    $ws.Range($chg.Cell).Value = $chg.Value
}

Write-Output ("Applied {0} cell updates" -f $changes.Count)
